# Update the cryptocurrency price/volume table (GitHub Actions refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D whose new text looks numeric ("307.73", "0.5281", ...)
# must be forced to Text format first, otherwise Excel silently converts
# the assigned string into a floating-point number (and "1.001"-style
# triple-grouped numbers like "27.315.45" would even mis-parse). Setting
# NumberFormat "@" then writing the literal string preserves the exact
# text; reapplying the "Normal" style afterwards keeps cell formatting
# identical to the original (no s="..." attribute) once saved.
$textCells = @("D5","D7","D8","D9","D10","D11","D12","D13","D14","D17","D18","D21","D23","D24","D25","D27","D28","D29","D30","D31","D32","D34","D35","D36","D37","D39","D40","D41","D42","D43","D44","D45","D47","D48","D50","D51")
foreach ($ref in $textCells) {
    $ws.Range($ref).NumberFormat = "@"
}

$ws.Range("D2").Value = '27.285.22'
$ws.Range("E2").Value = '  +0.34%  '

$ws.Range("D3").Value = '1.908.56'
$ws.Range("E3").Value = '  +0.34%  '

$ws.Range("E4").Value = '  +0.15%  '

$ws.Range("D5").Value = '307.73'

$ws.Range("E6").Value = '  +0.06%  '

$ws.Range("D7").Value = '0.5281'
$ws.Range("E7").Value = '  +1.32%  '

$ws.Range("D8").Value = '0.3821'
$ws.Range("E8").Value = '  +1.48%  '

$ws.Range("D9").Value = '0.07304'
$ws.Range("E9").Value = '  +0.39%  '

$ws.Range("D10").Value = '22.09'
$ws.Range("E10").Value = '  +4.49%  '

$ws.Range("D11").Value = '0.9022'
$ws.Range("E11").Value = '  -0.32%  '

$ws.Range("D12").Value = '0.08198'
$ws.Range("E12").Value = '  -1.04%  '

$ws.Range("D13").Value = '95.82'
$ws.Range("E13").Value = '  -0.93%  '

$ws.Range("D14").Value = '5.354'
$ws.Range("E14").Value = '  +1.19%  '

$ws.Range("E15").Value = '  +0.11%  '

$ws.Range("E16").Value = '  -0.35%  '

$ws.Range("D17").Value = '14.79'
$ws.Range("E17").Value = '  +1.51%  '

$ws.Range("B18").Value = 'Dai'
$ws.Range("C18").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D18").Value = '1.001'
$ws.Range("E18").Value = '  +0.15%  '

$ws.Range("B19").Value = 'WrappedEther'
$ws.Range("C19").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D19").Value = '1.296.42'
$ws.Range("E19").Value = '  -31.94%  '

$ws.Range("D20").Value = '27.315.45'
$ws.Range("E20").Value = '  +0.32%  '

$ws.Range("D21").Value = '5.068'
$ws.Range("E21").Value = '  -0.49%  '

$ws.Range("E22").Value = '  +1.57%  '

$ws.Range("D23").Value = '6.520'
$ws.Range("E23").Value = '  +1.25%  '

$ws.Range("D24").Value = '149.87'
$ws.Range("E24").Value = '  +2.37%  '

$ws.Range("D25").Value = '2.296'
$ws.Range("E25").Value = '  -1.12%  '

$ws.Range("E26").Value = '  +0.06%  '

$ws.Range("D27").Value = '1.738'
$ws.Range("E27").Value = '  -0.55%  '

$ws.Range("D28").Value = '116.45'
$ws.Range("E28").Value = '  +1.14%  '

$ws.Range("D29").Value = '4.827'
$ws.Range("E29").Value = '  -0.20%  '

$ws.Range("D30").Value = '4.820'
$ws.Range("E30").Value = '  -1.62%  '

$ws.Range("D31").Value = '0.09272'
$ws.Range("E31").Value = '  -0.02%  '

$ws.Range("D32").Value = '0.8367'
$ws.Range("E32").Value = '  +4.70%  '

$ws.Range("E33").Value = '  -0.19%  '

$ws.Range("D34").Value = '1.228'
$ws.Range("E34").Value = '  -1.36%  '

$ws.Range("D35").Value = '3.007'
$ws.Range("E35").Value = '  +2.18%  '

$ws.Range("D36").Value = '3.351'
$ws.Range("E36").Value = '  -2.04%  '

$ws.Range("D37").Value = '2.682'
$ws.Range("E37").Value = '  +3.50%  '

$ws.Range("E38").Value = '  +0.49%  '

$ws.Range("D39").Value = '0.02006'
$ws.Range("E39").Value = '  +0.21%  '

$ws.Range("D40").Value = '1.076'
$ws.Range("E40").Value = '  -0.07%  '

$ws.Range("D41").Value = '9.359'
$ws.Range("E41").Value = '  +3.72%  '

$ws.Range("D42").Value = '6.535'
$ws.Range("E42").Value = '  -0.84%  '

$ws.Range("D43").Value = '116.81'
$ws.Range("E43").Value = '  -0.28%  '

$ws.Range("D44").Value = '0.1522'
$ws.Range("E44").Value = '  +0.26%  '

$ws.Range("D45").Value = '0.4920'
$ws.Range("E45").Value = '  +1.11%  '

$ws.Range("E46").Value = '  +0.08%  '

$ws.Range("D47").Value = '10.17'
$ws.Range("E47").Value = '  +0.24%  '

$ws.Range("D48").Value = '1.637'
$ws.Range("E48").Value = '  +0.43%  '

$ws.Range("E49").Value = '  +2.98%  '

$ws.Range("D50").Value = '0.06185'
$ws.Range("E50").Value = '  +3.89%  '

$ws.Range("D51").Value = '63.75'
$ws.Range("E51").Value = '  -0.45%  '

# Restore default styling on the cells we text-formatted above.
foreach ($ref in $textCells) {
    $ws.Range($ref).Style = "Normal"
}
